# BAU RPS Data by Subregion.xlsx - apply the author's edits:
#   1. Rename sheet "BRPSDbS-RPS-percentage" -> "BRPSDbS-RPS-percentages"
#   2. Make that sheet the active sheet (it was "BRPSDbS-electricity-shares" before)
#   3. Move the active cell / selection on that sheet from F58 to C44

$wb = $excel.ActiveWorkbook

# 1. Rename the sheet.
$ws = $wb.Worksheets.Item("BRPSDbS-RPS-percentage")
$ws.Name = "BRPSDbS-RPS-percentages"

# 2. Switch to it so it becomes the workbook's active tab.
$ws.Activate() | Out-Null

# 3. Update the selection / active cell on the now-active sheet.
$ws.Range("C44").Select() | Out-Null
